$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4053036.5
$ws.Range("I86").Value = 5166.3335
$ws.Range("J86").Value = 5267397.5
$ws.Range("K86").Value = 5166.3335
$ws.Range("L86").Value = 5267397.5
$ws.Range("M86").Value = -4043.3335
$ws.Range("N86").Value = -5269643.5

$ws.Range("H89").Value = 4053036.5
$ws.Range("I89").Value = 5166.3335
$ws.Range("J89").Value = 5267397.5
$ws.Range("K89").Value = 25831.6675
$ws.Range("L89").Value = 26336987.5
$ws.Range("M89").Value = -20215.6675
$ws.Range("N89").Value = -26348219.5

$ws.Range("H106").Value = 4818.5
$ws.Range("I106").Value = 4818.5
$ws.Range("K106").Value = 4818.5
$ws.Range("M106").Value = -4187.5

$ws.Range("H132").Value = 2429.6
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 3000
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3950.8447
$ws.Range("I32").Value = 3576.0894
$ws.Range("K32").Value = 3576.0894
$ws.Range("M32").Value = -3289.0894

$ws.Range("H45").Value = 2417.0715
$ws.Range("I45").Value = 1918.2222
$ws.Range("J45").Value = 3315
$ws.Range("K45").Value = 1918.2222
$ws.Range("L45").Value = 3315
$ws.Range("M45").Value = -1541.2222
$ws.Range("N45").Value = -4069

$ws.Range("H74").Value = 1039.65
$ws.Range("I74").Value = 1039.65
$ws.Range("K74").Value = 1039.65
$ws.Range("M74").Value = -165.6500000000001

$ws.Range("H77").Value = 1039.65
$ws.Range("I77").Value = 1039.65
$ws.Range("K77").Value = 5198.25
$ws.Range("M77").Value = -830.25

$ws.Range("H132").Value = 3601.0667
$ws.Range("I132").Value = 3821.28
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 11463.84
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -8933.84
$ws.Range("N132").Value = -12560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 946097.5600000001
$ws.Range("I86").Value = 1215775.4
$ws.Range("K86").Value = 1215775.4
$ws.Range("M86").Value = -1214652.4

$ws.Range("H89").Value = 946097.5600000001
$ws.Range("I89").Value = 1215775.4
$ws.Range("K89").Value = 6078877
$ws.Range("M89").Value = -6073261

$ws.Range("H107").Value = 911550.5600000001
$ws.Range("I107").Value = 1630.625
$ws.Range("J107").Value = 3338003.8
$ws.Range("K107").Value = 1630.625
$ws.Range("L107").Value = 3338003.8
$ws.Range("M107").Value = 289.375
$ws.Range("N107").Value = -3341843.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30786
$ws.Range("I31").Value = 2590
$ws.Range("K31").Value = 2590
$ws.Range("M31").Value = -2295

$ws.Range("H34").Value = 30786
$ws.Range("I34").Value = 2590
$ws.Range("K34").Value = 2590
$ws.Range("M34").Value = -2388

$ws.Range("H58").Value = 5314.356
$ws.Range("I58").Value = 5067.024
$ws.Range("K58").Value = 5067.024
$ws.Range("M58").Value = -4864.024

$ws.Range("H122").Value = 3101.9546
$ws.Range("I122").Value = 2706.2727
$ws.Range("J122").Value = 3497.6365
$ws.Range("K122").Value = 8118.8181
$ws.Range("L122").Value = 10492.9095
$ws.Range("M122").Value = -5668.8181
$ws.Range("N122").Value = -15392.9095

$ws.Range("H132").Value = 2096.3333
$ws.Range("I132").Value = 2096.3333
$ws.Range("K132").Value = 6288.999899999999
$ws.Range("M132").Value = -3758.999899999999

$ws.Range("H134").Value = 373350.2
$ws.Range("I134").Value = 3094.423
$ws.Range("K134").Value = 9283.269
$ws.Range("M134").Value = -6748.269

$ws.Range("H136").Value = 5314.356
$ws.Range("I136").Value = 5067.024
$ws.Range("K136").Value = 15201.072
$ws.Range("M136").Value = -12651.072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 30395524
$ws.Range("J129").Value = 168839
$ws.Range("L129").Value = 506517
$ws.Range("N129").Value = -516517

$ws.Range("H136").Value = 5006
$ws.Range("I136").Value = 3757.5
$ws.Range("K136").Value = 11272.5
$ws.Range("M136").Value = -6172.5

$ws.Range("H138").Value = 2402.4614
$ws.Range("I138").Value = 2248.5454
$ws.Range("J138").Value = 3249
$ws.Range("K138").Value = 6745.6362
$ws.Range("L138").Value = 9747
$ws.Range("M138").Value = -1605.6362
$ws.Range("N138").Value = -20027

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1002423.4
$ws.Range("I80").Value = 772220.75
$ws.Range("J80").Value = 1429942.4
$ws.Range("K80").Value = 772220.75
$ws.Range("L80").Value = 1429942.4
$ws.Range("M80").Value = -771222.75
$ws.Range("N80").Value = -1431938.4

$ws.Range("H83").Value = 1002423.4
$ws.Range("I83").Value = 772220.75
$ws.Range("J83").Value = 1429942.4
$ws.Range("K83").Value = 3861103.75
$ws.Range("L83").Value = 7149712
$ws.Range("M83").Value = -3856111.75
$ws.Range("N83").Value = -7159696

$ws.Range("H132").Value = 39370.355
$ws.Range("J132").Value = 101844.9
$ws.Range("L132").Value = 305534.7
$ws.Range("N132").Value = -310594.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 500000740
$ws.Range("J16").Value = 1500
$ws.Range("L16").Value = 1500
$ws.Range("N16").Value = -1840

$ws.Range("H127").Value = 74489
$ws.Range("J127").Value = 74489
$ws.Range("L127").Value = 74489
$ws.Range("N127").Value = -84409

$ws.Range("H132").Value = 6444.143
$ws.Range("I132").Value = 4468.778
$ws.Range("K132").Value = 13406.334
$ws.Range("M132").Value = -10876.334

$ws.Range("H136").Value = 1117662.5
$ws.Range("I136").Value = 1673051.9
$ws.Range("K136").Value = 5019155.699999999
$ws.Range("M136").Value = -5016605.699999999

$ws.Range("H138").Value = 93326.664
$ws.Range("J138").Value = 93326.664
$ws.Range("L138").Value = 93326.664
$ws.Range("N138").Value = -103606.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1609.75
$ws.Range("I81").Value = 1434.6666
$ws.Range("K81").Value = 2869.3332
$ws.Range("M81").Value = -1808.3332

$ws.Range("H84").Value = 1609.75
$ws.Range("I84").Value = 1434.6666
$ws.Range("K84").Value = 14346.666
$ws.Range("M84").Value = -9042.666000000001

$ws.Range("H107").Value = 605.63635
$ws.Range("I107").Value = 661.88
$ws.Range("J107").Value = 429.875
$ws.Range("K107").Value = 1985.64
$ws.Range("L107").Value = 1289.625
$ws.Range("M107").Value = -65.63999999999987
$ws.Range("N107").Value = -5129.625

$ws.Range("H122").Value = 62502300
$ws.Range("I122").Value = 83334790
$ws.Range("K122").Value = 250004370
$ws.Range("M122").Value = -250001920

$ws.Range("H132").Value = 24096.936
$ws.Range("I132").Value = 1624.1
$ws.Range("K132").Value = 4872.299999999999
$ws.Range("M132").Value = -2342.299999999999

$ws.Range("H136").Value = 9345873
$ws.Range("I136").Value = 14324278
$ws.Range("J136").Value = 154969.84
$ws.Range("K136").Value = 42972834
$ws.Range("L136").Value = 464909.52
$ws.Range("M136").Value = -42970284
$ws.Range("N136").Value = -470009.52
